$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Strike through the whole "Agregar botones Edit y Delete..." bullet item
#    (the task is done, so the author marked it as struck-through).
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "Agregar botones Edit y Delete en las tarjetas Experience, Education y Projects (A implementar)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # Use the whole containing paragraph (includes the paragraph mark) so
    # the strike-through also lands on the paragraph-mark run properties,
    # matching how Word stores "select paragraph, strike it" edits.
    $para = $rng.Paragraphs(1)
    $para.Range.Font.StrikeThrough = $true
}

# ---------------------------------------------------------------------------
# 2) Insert " para esto" right before ") (A testear)" in the
#    "Hacer vertical Contact..." bullet item.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "cols) (A testear)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $suffixLen = ") (A testear)".Length
    $target = $d.Range($rng2.End - $suffixLen, $rng2.End)
    $target.Text = " para esto" + $target.Text
}

# ---------------------------------------------------------------------------
# 3) Append a new sentence at the end of the "16/06/22" log entry paragraph,
#    right after "...se acomodan al tamaño de la pantalla."
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute(
    "se acomodan al tamaño de la pantalla.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $insertAt = $rng3.End
    $newText = " Y agregue los botones de edición en los componentes que lo necesitaban y elimine botones en donde sobraban. La mayoría de los botones tiene posición absoluta para que no interfieran con otros elementos."
    $rng3.InsertAfter($newText)

    # Match the formatting of the surrounding paragraph text (Courier New,
    # black, eastAsia font Times New Roman) by driving the replacement
    # through Find/Replace formatting, which is what actually serializes
    # the full run-fonts set (ascii/eastAsia/hAnsi/cs) in this runtime.
    $newRange = $d.Range($insertAt, $insertAt + $newText.Length)
    $newRange.Find.ClearFormatting()
    $newRange.Find.Replacement.ClearFormatting()
    $newRange.Find.Replacement.Font.Name = "Courier New"
    $newRange.Find.Replacement.Font.NameFarEast = "Times New Roman"
    $newRange.Find.Replacement.Font.NameOther = "Courier New"
    $newRange.Find.Replacement.Font.NameBi = "Courier New"
    $newRange.Find.Replacement.Font.Color = 0
    $newRange.Find.Execute($newText, $true, $false, $false, $false, $false, $true, 1, $true, $newText, 2)
}
